$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header-adjacent data rows 2-7 to reflect revised NATMI LR-pair output
# (adds ECs sending/target cluster combinations per Dr Hou advice)

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Wnt1"
$ws.Range("C2").Value = "Fzd2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.027123
$ws.Range("H2").Value = 0.081369
$ws.Range("I2").Value = 0.07131444737854614
$ws.Range("J2").Value = 0.07131444737854616
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.07629866666666667
$ws.Range("N2").Value = 0.228896
$ws.Range("O2").Value = 0.004108848954870246
$ws.Range("P2").Value = 0.004108848954870246
$ws.Range("Q2").Value = 0.002069448736
$ws.Range("R2").Value = 0.018625038624
$ws.Range("S2").Value = 0.0002930202925784885
$ws.Range("T2").Value = 0.0002930202925784885

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Wnt1"
$ws.Range("C3").Value = "Fzd2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.027123
$ws.Range("H3").Value = 0.081369
$ws.Range("I3").Value = 0.07131444737854614
$ws.Range("J3").Value = 0.07131444737854616
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 12.99468133333333
$ws.Range("N3").Value = 38.984044
$ws.Range("O3").Value = 0.6997918200668237
$ws.Range("P3").Value = 0.6997918200668237
$ws.Range("Q3").Value = 0.3524547418039999
$ws.Range("R3").Value = 3.172092676236
$ws.Range("S3").Value = 0.04990526692809253
$ws.Range("T3").Value = 0.04990526692809254

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Wnt1"
$ws.Range("C4").Value = "Fzd2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.027123
$ws.Range("H4").Value = 0.081369
$ws.Range("I4").Value = 0.07131444737854614
$ws.Range("J4").Value = 0.07131444737854616
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 5.498373000000001
$ws.Range("N4").Value = 16.495119
$ws.Range("O4").Value = 0.2960993309783061
$ws.Range("P4").Value = 0.2960993309783061
$ws.Range("Q4").Value = 0.149132370879
$ws.Range("R4").Value = 1.342191337911
$ws.Range("S4").Value = 0.02111616015787513
$ws.Range("T4").Value = 0.02111616015787513

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Wnt1"
$ws.Range("C5").Value = "Fzd2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.3532066666666667
$ws.Range("H5").Value = 1.05962
$ws.Range("I5").Value = 0.9286855526214538
$ws.Range("J5").Value = 0.9286855526214538
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.07629866666666667
$ws.Range("N5").Value = 0.228896
$ws.Range("O5").Value = 0.004108848954870246
$ws.Range("P5").Value = 0.004108848954870246
$ws.Range("Q5").Value = 0.02694919772444445
$ws.Range("R5").Value = 0.24254277952
$ws.Range("S5").Value = 0.003815828662291758
$ws.Range("T5").Value = 0.003815828662291758

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Wnt1"
$ws.Range("C6").Value = "Fzd2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.3532066666666667
$ws.Range("H6").Value = 1.05962
$ws.Range("I6").Value = 0.9286855526214538
$ws.Range("J6").Value = 0.9286855526214538
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 12.99468133333333
$ws.Range("N6").Value = 38.984044
$ws.Range("O6").Value = 0.6997918200668237
$ws.Range("P6").Value = 0.6997918200668237
$ws.Range("Q6").Value = 4.589808078142222
$ws.Range("R6").Value = 41.30827270328
$ws.Range("S6").Value = 0.6498865531387311
$ws.Range("T6").Value = 0.6498865531387311

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Wnt1"
$ws.Range("C7").Value = "Fzd2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.3532066666666667
$ws.Range("H7").Value = 1.05962
$ws.Range("I7").Value = 0.9286855526214538
$ws.Range("J7").Value = 0.9286855526214538
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.498373000000001
$ws.Range("N7").Value = 16.495119
$ws.Range("O7").Value = 0.2960993309783061
$ws.Range("P7").Value = 0.2960993309783061
$ws.Range("Q7").Value = 1.94206199942
$ws.Range("R7").Value = 17.47855799478
$ws.Range("S7").Value = 0.2749831708204309
$ws.Range("T7").Value = 0.2749831708204309

